# Update kaspa_buys.xlsx after running on 2025-05-15
# Appends a new purchase record as row 24 (Date, Coins, Price, Cost).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 24

# Column A holds the date as literal text (e.g. "05/15/2025"), matching the
# other recently-appended rows in this sheet. Force the cell to Text format
# before assigning the value so Excel doesn't auto-convert the string into a
# date serial number, then strip the temporary formatting back off so the
# cell is left with no explicit style (same as the existing text-date rows).
$dateCell = $ws.Range("A" + $newRow)
$dateCell.NumberFormat = "@"
$dateCell.Value = "05/15/2025"
$dateCell.ClearFormats()

$ws.Range("B" + $newRow).Value = 414.3919999999998
$ws.Range("C" + $newRow).Value = 0.120658699974903
$ws.Range("D" + $newRow).Value = 50
